$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = '139 Highett St Apartment Complex Richmond'
$ws.Cells.Item(2, 2).Value = 11
$ws.Cells.Item(3, 1).Value = '3175 The Bays Aged Care Facility Hastings'
$ws.Cells.Item(3, 2).Value = 16
$ws.Cells.Item(4, 1).Value = '3600 Belvedere Age Care Noble Park Outbreak'
$ws.Cells.Item(4, 2).Value = 21
$ws.Cells.Item(5, 1).Value = '3612 BlueCross Glengowrie Outbreak'
$ws.Cells.Item(5, 2).Value = 30
$ws.Cells.Item(6, 1).Value = '3684 Homestyle Aged Care Langford Grange Cranbourne East Outbreak'
$ws.Cells.Item(6, 2).Value = 20
$ws.Cells.Item(7, 1).Value = '4075 Ferndale Gardens Aged Care Services Bayswater North Outbreak'
$ws.Cells.Item(7, 2).Value = 16
$ws.Cells.Item(8, 1).Value = 'Australian Lamb Colac East'
$ws.Cells.Item(8, 2).Value = 13
$ws.Cells.Item(9, 1).Value = 'Bread Solutions Braeside Outbreak'
$ws.Cells.Item(9, 2).Value = 19
$ws.Cells.Item(10, 1).Value = 'CS Square Caroline Springs Outbreak'
$ws.Cells.Item(10, 2).Value = 17
$ws.Cells.Item(11, 1).Value = 'Cedar Meats Australia Brooklyn Outbreak'
$ws.Cells.Item(11, 2).Value = 11
$ws.Cells.Item(12, 1).Value = 'Child''s Play Early Learning Centre Tarneit'
$ws.Cells.Item(12, 2).Value = 11
$ws.Cells.Item(13, 1).Value = 'Embracia Aged Care Reservoir Outbreak'
$ws.Cells.Item(13, 2).Value = 22
$ws.Cells.Item(14, 1).Value = 'Guardian Childcare Caulfield Outbreak'
$ws.Cells.Item(14, 2).Value = 17
$ws.Cells.Item(15, 1).Value = 'Hello Fresh Warehouse Ravenhall'
$ws.Cells.Item(15, 2).Value = 13
$ws.Cells.Item(16, 1).Value = 'Inghams Enterprise Somerville Outbreak'
$ws.Cells.Item(16, 2).Value = 15
$ws.Cells.Item(17, 1).Value = 'Kool Kidz Childcare Narre Warren'
$ws.Cells.Item(17, 2).Value = 10
$ws.Cells.Item(18, 1).Value = 'Lantmannen Unibake Australia Mordialloc'
$ws.Cells.Item(18, 2).Value = 25
$ws.Cells.Item(19, 1).Value = 'Launch Housing City Edge Crisis Accommodation South Melbourne'
$ws.Cells.Item(19, 2).Value = 12
$ws.Cells.Item(20, 1).Value = 'Nido Early School Ascot Vale'
$ws.Cells.Item(20, 2).Value = 11
$ws.Cells.Item(21, 1).Value = 'Northern Health Northern Hospital Epping Emergency Department Tier 1B'
$ws.Cells.Item(21, 2).Value = 42
$ws.Cells.Item(22, 1).Value = 'Northern Health The Northern Hospital Epping'
$ws.Cells.Item(22, 2).Value = 15
$ws.Cells.Item(23, 1).Value = 'Oceania Meat Processors Laverton North Outbreak'
$ws.Cells.Item(23, 2).Value = 16
$ws.Cells.Item(24, 1).Value = 'Pick It Up Fitness Mulgrave Outbreak'
$ws.Cells.Item(24, 2).Value = 11
$ws.Cells.Item(25, 1).Value = 'Robin Hood Inn Drouin West Outbreak'
$ws.Cells.Item(25, 2).Value = 41
$ws.Cells.Item(26, 1).Value = 'Social Gathering Warrnambool 28 Sep Outbreak'
$ws.Cells.Item(26, 2).Value = 13
$ws.Cells.Item(27, 1).Value = 'St Vincents Hospital Emergency Department Melbourne'
$ws.Cells.Item(27, 2).Value = 42
$ws.Cells.Item(28, 1).Value = 'Target Distribution Centre Truganina Outbreak'
$ws.Cells.Item(28, 2).Value = 20
$ws.Cells.Item(29, 1).Value = 'The Royal Children''s Hospital Melbourne Emergency Department Parkville Tier 1B'
$ws.Cells.Item(29, 2).Value = 12
$ws.Cells.Item(30, 1).Value = 'The Toolshed Bar Private Event Noojee'
$ws.Cells.Item(30, 2).Value = 17
$ws.Cells.Item(31, 1).Value = 'Turosi Breakwater'
$ws.Cells.Item(31, 2).Value = 10
$ws.Cells.Item(32, 1).Value = 'Visy Recycling Springvale'
$ws.Cells.Item(32, 2).Value = 29
$ws.Cells.Item(33, 1).Value = 'Werribee Mercy Hospital Emergency Department'
$ws.Cells.Item(33, 2).Value = 24
$ws.Cells.Item(34, 1).Value = 'Western Health Sunshine Hospital Emergency Department'
$ws.Cells.Item(34, 2).Value = 21
# Remove the now-extra 35th row (table shrank from 34 data rows to 33 data rows)
$ws.Rows.Item(35).ClearContents()
